# -----------------------------------------------------------------------
# This script reproduces the commit "update main file to add more
# experimentation": it updates a handful of existing values on sheets
# "cost" (sheet2), "availability" (sheet3), "reputation" (sheet4) and
# "time" (sheet5), rewrites the ratio formulas on "reputation" to use
# straightforward division, appends a new data row (dataset size 3000)
# to every sheet, removes the leftover "min is prefered" note from the
# "reputation" sheet, adds print setup to the "all" sheet, and finally
# makes "reputation" the active tab.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsAll  = $wb.Worksheets.Item("all")
$wsCost = $wb.Worksheets.Item("cost")
$wsAvail = $wb.Worksheets.Item("availability")
$wsRep  = $wb.Worksheets.Item("reputation")
$wsTime = $wb.Worksheets.Item("time")

# ------------------------------------------------------------------
# Sheet "all" (sheet1)
# ------------------------------------------------------------------
$wsAll.Range("F3").Value = 4.1067299840713103

$wsAll.Range("A11").Value = 3000
$wsAll.Range("B11").Value = 5.7642854309953417
$wsAll.Range("C11").Value = 5.2548878189512456
$wsAll.Range("D11").Value = 5.0099045352888938
$wsAll.Range("E11").Value = 5.3861770710582206
$wsAll.Range("F11").Value = 4.7039457738836461

$wsAll.PageSetup.PaperSize = 9
$wsAll.PageSetup.Orientation = 1

$wsAll.Range("L13").Select()

# ------------------------------------------------------------------
# Sheet "cost" (sheet2)
# ------------------------------------------------------------------
$wsCost.Range("F2").Value = 151.012973833345
$wsCost.Range("F3").Value = 273.27702853708797
$wsCost.Range("F4").Value = 150.575649185792
$wsCost.Range("F5").Value = 138.052510748797
$wsCost.Range("F6").Value = 123.327246835701
$wsCost.Range("F7").Value = 108.205892984267
$wsCost.Range("F8").Value = 134.21483768904699

$wsCost.Range("A11").Value = 3000
$wsCost.Range("B11").Value = 656.73414490522839
$wsCost.Range("C11").Value = 341.31962901642339
$wsCost.Range("D11").Value = 121.16006895004389
$wsCost.Range("E11").Value = 510.74240433444038
$wsCost.Range("F11").Value = 115.660598947483

$wsCost.Range("F11").Select()

# ------------------------------------------------------------------
# Sheet "availability" (sheet3)
# ------------------------------------------------------------------
$wsAvail.Range("F2").Value = 1.07036935103784
$wsAvail.Range("F3").Value = 1.15542158350902
$wsAvail.Range("F4").Value = 1.1286566038577901
$wsAvail.Range("F6").Value = 1.08733345397738
$wsAvail.Range("F7").Value = 1.1403687912342499
$wsAvail.Range("F8").Value = 1.08470513579211
$wsAvail.Range("D9").Value = 1.3453094821976299
$wsAvail.Range("F9").Value = 1.18479162356969
$wsAvail.Range("D10").Value = 1.5658986209808099

$wsAvail.Range("A11").Value = 3000
$wsAvail.Range("B11").Value = 4.5085080125366428
$wsAvail.Range("C11").Value = 1.7582376736068701
$wsAvail.Range("D11").Value = 1.159012842724342
$wsAvail.Range("E11").Value = 2.511749782748304
$wsAvail.Range("F11").Value = 1.0226802392191929

$wsAvail.Range("F11").Select()

# ------------------------------------------------------------------
# Sheet "reputation" (sheet4)
# ------------------------------------------------------------------
$wsRep.Range("F2").Value = 314.91334617713602
$wsRep.Range("F3").Value = 633.73508556486604
$wsRep.Range("F4").Value = 854.85438633621004
$wsRep.Range("F6").Value = 1575.91698096954
$wsRep.Range("F7").Value = 1727.3102195239301
$wsRep.Range("F8").Value = 2109.1885946452999
$wsRep.Range("F10").Value = 3280.4239198710602

# Convert the ratio formulas in H:L from "(1/x)*y" style to plain "y/x"
# for every existing row.
for ($r = 2; $r -le 10; $r++) {
    $wsRep.Range("H$r").Formula = "=B$r/A$r"
    $wsRep.Range("I$r").Formula = "=C$r/A$r"
    $wsRep.Range("J$r").Formula = "=D$r/A$r"
    $wsRep.Range("K$r").Formula = "=E$r/A$r"
    $wsRep.Range("L$r").Formula = "=F$r/A$r"
}

# Drop the old "min is prefered" note that used to sit below the table.
$wsRep.Range("E13").ClearContents()

# Append the new dataset-size-3000 row together with its ratio formulas.
$wsRep.Range("A11").Value = 3000
$wsRep.Range("B11").Value = 7175.1406025074339
$wsRep.Range("C11").Value = 5521.8215821257081
$wsRep.Range("D11").Value = 3225.8936606957368
$wsRep.Range("E11").Value = 28145.988185038241
$wsRep.Range("F11").Value = 3848.9102960985902

$wsRep.Range("H11").Formula = "=B11/A11"
$wsRep.Range("I11").Formula = "=C11/A11"
$wsRep.Range("J11").Formula = "=D11/A11"
$wsRep.Range("K11").Formula = "=E11/A11"
$wsRep.Range("L11").Formula = "=F11/A11"

# ------------------------------------------------------------------
# Sheet "time" (sheet5)
# ------------------------------------------------------------------
$wsTime.Range("C2").Value = 1.17744144637619
$wsTime.Range("D2").Value = 1.55288157075202
$wsTime.Range("C3").Value = 1.2346927306748501
$wsTime.Range("D3").Value = 0.95125054694352495
$wsTime.Range("C4").Value = 1.1785839173452901
$wsTime.Range("D4").Value = 1.5523891728982
$wsTime.Range("C5").Value = 1.1961687775343399
$wsTime.Range("D5").Value = 1.5648899337217499
$wsTime.Range("C6").Value = 1.2625971669060301
$wsTime.Range("F6").Value = 1.0052914878756301
$wsTime.Range("D7").Value = 0.98536048510236696
$wsTime.Range("D8").Value = 1.1207602820185101
$wsTime.Range("C9").Value = 1.32542360138366
$wsTime.Range("D9").Value = 1.2824878127011601
$wsTime.Range("C10").Value = 1.56761132158565
$wsTime.Range("D10").Value = 1.15610889905773

$wsTime.Range("A11").Value = 3000
$wsTime.Range("B11").Value = 3.7990695122187308
$wsTime.Range("C11").Value = 1.529220721322408
$wsTime.Range("D11").Value = 1.57437225712201
$wsTime.Range("E11").Value = 1.7108662717181
$wsTime.Range("F11").Value = 1.5336274523315649

$wsTime.Range("F11").Select()

# Activate "reputation" last so it ends up as the workbook's active tab,
# matching the saved view state of the edited workbook.
$wsRep.Activate()
$wsRep.Range("H3").Select()

Write-Host "edit applied"
